$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column B for "Variable/Constant" ---
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "Variable/Constant"

# --- Insert blank rows that will hold the new "Constant" rows ---
# Work bottom-to-top (by current row number) so earlier inserts don't
# disturb the row numbers used by later inserts.
$ws.Rows("4:4").Insert()
$ws.Rows("3:3").Insert()
$ws.Rows("2:2").Insert()

# Newly-inserted rows copy the formatting (style + content placeholders)
# of the row above them. Strip that back out so only the cells we
# actually populate end up in the saved XML.
$ws.Rows("2:2").ClearFormats()
$ws.Rows("4:4").ClearFormats()
$ws.Rows("6:6").ClearFormats()
$ws.Range("D2:J2").ClearContents()
$ws.Range("D4:J4").ClearContents()
$ws.Range("D6:J6").ClearContents()

# At this point the sheet looks like:
#   Row1 = header
#   Row2 = blank (new "Constant" row)
#   Row3 = old row2 data (202 / Top)
#   Row4 = blank (new "Constant" row)
#   Row5 = old row3 data (231 / Top)
#   Row6 = blank (new "Constant" row)
#   Row7 = old row4 data (242 / Top)
# Rows 8-10 are brand new (appended below the existing data).

# --- Label every row's Region / Protein / Variable-Constant columns ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 3).Value = "MHC A3"
    if ($r % 2 -eq 0) {
        $ws.Cells.Item($r, 2).Value = "Constant"
    } else {
        $ws.Cells.Item($r, 2).Value = "Variable"
    }
}

# --- Fill in the "Constant" rows (2, 4, 6, 8, 10) with their placeholder residue ---
$ws.Cells.Item(2, 9).Value = "aaaa"
$ws.Cells.Item(4, 9).Value = "bbbb"
$ws.Cells.Item(6, 9).Value = "cccc"
$ws.Cells.Item(8, 9).Value = "dddd"
$ws.Cells.Item(10, 9).Value = "eeee"

# --- Fill in the new 9th data row (row 9) - new "Variable" entry ---
$ws.Cells.Item(9, 4).Value = 246
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = "Top"
$ws.Cells.Item(9, 7).Value = "GTTGTCGTTCCATCT"
$ws.Cells.Item(9, 8).Value = "CAAGAACAGCGTTAT"
$ws.Cells.Item(9, 9).Value = "G"
$ws.Cells.Item(9, 10).Value = "GT"

# --- Page setup (portrait orientation, as recorded in the saved workbook) ---
$ws.PageSetup.Orientation = 1

# --- Selection / view tweaks ---
$ws.Range("B11").Select()
